$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in columns P, Q, R for rows 4 and 5 ---
$ws.Range("P4").Value = 13.7
$ws.Range("Q4").Value = 13.1
$ws.Range("R4").Value = 11.8

$ws.Range("P5").Value = 13.6
$ws.Range("Q5").Value = 12.5
$ws.Range("R5").Value = 13.5

# --- Add new column S, mirroring the formatting of column R for each row ---
$ws.Range("R2").Copy()
$ws.Range("S2").PasteSpecial(-4122)

$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial(-4122)
$ws.Range("S3").Value = 2022

$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 13.6

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 20

$excel.CutCopyMode = 0

# --- Update selection to match the new active cell ---
[void]$ws.Range("S2").Select()
